$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Block 1 (ChiNext 50 / 创业板50（159949）): update the 5F row value pair
# B9 holds the "中枢1底/中枢1顶" style value for the 5F level; its execution-status
# marker (C9) and paired range (D9) are cleared since the cycle hasn't been confirmed yet.
$ws.Range("B9").Value = "1.09/1.107"
$ws.Range("C9").ClearContents()
$ws.Range("D9").ClearContents()

# Block 2 (300ETF（510300）): same kind of update for its 5F row
$ws.Range("B18").Value = "4.761/4.804"
$ws.Range("C18").ClearContents()
$ws.Range("D18").ClearContents()

# Update the active selection recorded in the sheet view
$ws.Range("C22").Select()
